{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// --- Collapse \"Test \" + \"git\" + \" 1\" (split across runs with spell-check\n// markers) into a single clean run \"Test git 1\" in the first paragraph. ---\nconst firstParagraph = paragraphs.items[0];\nfirstParagraph.getRange().insertText(\"Test git 1\", \"Replace\");\n\n// --- Append new paragraphs \"2\"..\"9\",\"0\" after the first paragraph. ---\nconst values = [\"2\", \"3\", \"4\", \"5\", \"6\", \"7\", \"8\", \"9\", \"0\"];\nlet anchor = firstParagraph;\nfor (const value of values) {\n  anchor = anchor.insertParagraph(value, \"After\");\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Collapse \"Test \" + \"git\" + \" 1\" (split across runs with spell-check\n# markers) into a single clean run \"Test git 1\" in the first paragraph. ---\n$firstRange = $d.Paragraphs(1).Range\n$firstRange.MoveEnd(1, -1) | Out-Null   # exclude the paragraph mark\n$firstRange.Text = \"\"                   # clear so the re-type below always registers\n$firstRange.Text = \"Test git 1\"\n\n# --- Append new paragraphs \"2\"..\"9\",\"0\" after the first paragraph. ---\n$values = @(\"2\", \"3\", \"4\", \"5\", \"6\", \"7\", \"8\", \"9\", \"0\")\nforeach ($v in $values) {\n    $lastParaRange = $d.Paragraphs.Last.Range\n    $lastParaRange.InsertParagraphAfter()\n    $newRange = $d.Paragraphs.Last.Range\n    $newRange.MoveEnd(1, -1) | Out-Null\n    $newRange.Text = $v\n}\n"}
